{"js": "const REPLACEMENTS = [\n  [\"2023-05-23 Tuesday\", \"2023-05-24 Wednesday\"],\n  [\"83-31=\", \"10+71=\"],\n  [\"68-19=\", \"35-5=\"],\n  [\"8+66=\", \"56-54=\"],\n  [\"86-25=\", \"38-23=\"],\n  [\"56+0=\", \"42-17=\"],\n  [\"42+47=\", \"10+47=\"],\n  [\"85+7=\", \"67-8=\"],\n  [\"67-35=\", \"54+7=\"],\n  [\"60-1=\", \"29+25=\"],\n  [\"93-26=\", \"78-52=\"],\n  [\"93-41=\", \"67+13=\"],\n  [\"26-21=\", \"87-40=\"],\n  [\"35+36=\", \"25+15=\"],\n  [\"42+19=\", \"47+31=\"],\n  [\"63-10=\", \"47-21=\"],\n  [\"94-3=\", \"87-81=\"],\n  [\"65-51=\", \"3+87=\"],\n  [\"76-74=\", \"77-40=\"],\n  [\"55-3=\", \"70-42=\"],\n  [\"41-20=\", \"41-4=\"],\n  [\"23+2=\", \"71-65=\"],\n  [\"31+43=\", \"75-15=\"],\n  [\"0+96=\", \"73-16=\"],\n  [\"5+4=\", \"74-27=\"],\n  [\"86-53=\", \"63-24=\"],\n  [\"41+37=\", \"96-37=\"],\n  [\"16+6=\", \"58-6=\"],\n  [\"61-19=\", \"90-5=\"],\n  [\"25+44=\", \"97-27=\"],\n  [\"42-26=\", \"80-18=\"],\n  [\"82-64=\", \"74-66=\"],\n  [\"31+64=\", \"54-26=\"],\n  [\"66+24=\", \"70-49=\"],\n  [\"61-52=\", \"10+38=\"],\n  [\"52-48=\", \"79+19=\"],\n  [\"14-1=\", \"25-13=\"],\n  [\"44+49=\", \"46-21=\"],\n  [\"18+31=\", \"88-73=\"],\n  [\"12+71=\", \"89-41=\"],\n  [\"94-55=\", \"81-15=\"],\n  [\"93-48=\", \"93-70=\"],\n  [\"72+6=\", \"65+9=\"],\n  [\"94-38=\", \"49+35=\"],\n  [\"6+61=\", \"90-88=\"],\n  [\"60+17=\", \"55+20=\"],\n  [\"80+10=\", \"41+39=\"],\n  [\"0+47=\", \"22+53=\"],\n  [\"1+55=\", \"16-13=\"],\n  [\"71+25=\", \"95-42=\"],\n  [\"35-23=\", \"71+18=\"],\n  [\"69-25=\", \"5+27=\"],\n  [\"12+87=\", \"2+39=\"],\n  [\"27+23=\", \"70-19=\"],\n  [\"0+92=\", \"26-15=\"],\n  [\"35-24=\", \"55+4=\"],\n  [\"24+37=\", \"97-8=\"],\n  [\"21+8=\", \"48-22=\"],\n  [\"89+3=\", \"86-16=\"],\n  [\"8+70=\", \"14+76=\"],\n  [\"25+70=\", \"46+39=\"],\n  [\"42-12=\", \"30+56=\"],\n  [\"76-4=\", \"98-48=\"],\n  [\"36+10=\", \"16+15=\"],\n  [\"3+10=\", \"57+42=\"],\n  [\"26+6=\", \"52+19=\"],\n  [\"70-69=\", \"44+23=\"],\n  [\"72+15=\", \"61+23=\"],\n  [\"11+54=\", \"68+17=\"],\n  [\"51+25=\", \"53-28=\"],\n  [\"76+7=\", \"21+20=\"],\n  [\"76-43=\", \"78-62=\"],\n  [\"38-31=\", \"12+18=\"],\n  [\"68-65=\", \"66+17=\"],\n  [\"9+7=\", \"37-23=\"],\n  [\"86-47=\", \"64+1=\"],\n  [\"37+5=\", \"86+5=\"],\n  [\"90-49=\", \"90-78=\"],\n  [\"75-74=\", \"90-47=\"],\n  [\"93-23=\", \"68+19=\"],\n  [\"58-12=\", \"98-1=\"],\n  [\"26+17=\", \"11+38=\"],\n  [\"89-87=\", \"63-1=\"],\n  [\"91-81=\", \"78-49=\"],\n  [\"20+41=\", \"1+30=\"],\n  [\"94-42=\", \"73-69=\"],\n  [\"17+44=\", \"37+60=\"],\n  [\"25+1=\", \"35+33=\"],\n  [\"78-17=\", \"39-39=\"],\n  [\"70-52=\", \"79+17=\"],\n  [\"81-49=\", \"12+5=\"],\n  [\"3+73=\", \"92-71=\"],\n  [\"96-85=\", \"13+80=\"],\n  [\"54-48=\", \"1+59=\"],\n  [\"31+56=\", \"87-55=\"],\n  [\"23+28=\", \"99-20=\"],\n  [\"27+33=\", \"44+43=\"],\n  [\"61+31=\", \"66-9=\"],\n  [\"18+18=\", \"67-64=\"],\n  [\"0+66=\", \"91-59=\"],\n  [\"90-71=\", \"32+3=\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet applied = 0;\nlet mismatches = [];\nfor (let i = 0; i < items.length && i < REPLACEMENTS.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  if (para.text === oldText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n    applied++;\n  } else {\n    mismatches.push({ index: i, expected: oldText, actual: para.text });\n  }\n}\nawait context.sync();\n\nreturn JSON.stringify({ applied: applied, total: REPLACEMENTS.length, mismatches: mismatches });\n", "ps1": "$Replacements = @(\n    @(\"2023-05-23 Tuesday\", \"2023-05-24 Wednesday\"),\n    @(\"83-31=\", \"10+71=\"),\n    @(\"68-19=\", \"35-5=\"),\n    @(\"8+66=\", \"56-54=\"),\n    @(\"86-25=\", \"38-23=\"),\n    @(\"56+0=\", \"42-17=\"),\n    @(\"42+47=\", \"10+47=\"),\n    @(\"85+7=\", \"67-8=\"),\n    @(\"67-35=\", \"54+7=\"),\n    @(\"60-1=\", \"29+25=\"),\n    @(\"93-26=\", \"78-52=\"),\n    @(\"93-41=\", \"67+13=\"),\n    @(\"26-21=\", \"87-40=\"),\n    @(\"35+36=\", \"25+15=\"),\n    @(\"42+19=\", \"47+31=\"),\n    @(\"63-10=\", \"47-21=\"),\n    @(\"94-3=\", \"87-81=\"),\n    @(\"65-51=\", \"3+87=\"),\n    @(\"76-74=\", \"77-40=\"),\n    @(\"55-3=\", \"70-42=\"),\n    @(\"41-20=\", \"41-4=\"),\n    @(\"23+2=\", \"71-65=\"),\n    @(\"31+43=\", \"75-15=\"),\n    @(\"0+96=\", \"73-16=\"),\n    @(\"5+4=\", \"74-27=\"),\n    @(\"86-53=\", \"63-24=\"),\n    @(\"41+37=\", \"96-37=\"),\n    @(\"16+6=\", \"58-6=\"),\n    @(\"61-19=\", \"90-5=\"),\n    @(\"25+44=\", \"97-27=\"),\n    @(\"42-26=\", \"80-18=\"),\n    @(\"82-64=\", \"74-66=\"),\n    @(\"31+64=\", \"54-26=\"),\n    @(\"66+24=\", \"70-49=\"),\n    @(\"61-52=\", \"10+38=\"),\n    @(\"52-48=\", \"79+19=\"),\n    @(\"14-1=\", \"25-13=\"),\n    @(\"44+49=\", \"46-21=\"),\n    @(\"18+31=\", \"88-73=\"),\n    @(\"12+71=\", \"89-41=\"),\n    @(\"94-55=\", \"81-15=\"),\n    @(\"93-48=\", \"93-70=\"),\n    @(\"72+6=\", \"65+9=\"),\n    @(\"94-38=\", \"49+35=\"),\n    @(\"6+61=\", \"90-88=\"),\n    @(\"60+17=\", \"55+20=\"),\n    @(\"80+10=\", \"41+39=\"),\n    @(\"0+47=\", \"22+53=\"),\n    @(\"1+55=\", \"16-13=\"),\n    @(\"71+25=\", \"95-42=\"),\n    @(\"35-23=\", \"71+18=\"),\n    @(\"69-25=\", \"5+27=\"),\n    @(\"12+87=\", \"2+39=\"),\n    @(\"27+23=\", \"70-19=\"),\n    @(\"0+92=\", \"26-15=\"),\n    @(\"35-24=\", \"55+4=\"),\n    @(\"24+37=\", \"97-8=\"),\n    @(\"21+8=\", \"48-22=\"),\n    @(\"89+3=\", \"86-16=\"),\n    @(\"8+70=\", \"14+76=\"),\n    @(\"25+70=\", \"46+39=\"),\n    @(\"42-12=\", \"30+56=\"),\n    @(\"76-4=\", \"98-48=\"),\n    @(\"36+10=\", \"16+15=\"),\n    @(\"3+10=\", \"57+42=\"),\n    @(\"26+6=\", \"52+19=\"),\n    @(\"70-69=\", \"44+23=\"),\n    @(\"72+15=\", \"61+23=\"),\n    @(\"11+54=\", \"68+17=\"),\n    @(\"51+25=\", \"53-28=\"),\n    @(\"76+7=\", \"21+20=\"),\n    @(\"76-43=\", \"78-62=\"),\n    @(\"38-31=\", \"12+18=\"),\n    @(\"68-65=\", \"66+17=\"),\n    @(\"9+7=\", \"37-23=\"),\n    @(\"86-47=\", \"64+1=\"),\n    @(\"37+5=\", \"86+5=\"),\n    @(\"90-49=\", \"90-78=\"),\n    @(\"75-74=\", \"90-47=\"),\n    @(\"93-23=\", \"68+19=\"),\n    @(\"58-12=\", \"98-1=\"),\n    @(\"26+17=\", \"11+38=\"),\n    @(\"89-87=\", \"63-1=\"),\n    @(\"91-81=\", \"78-49=\"),\n    @(\"20+41=\", \"1+30=\"),\n    @(\"94-42=\", \"73-69=\"),\n    @(\"17+44=\", \"37+60=\"),\n    @(\"25+1=\", \"35+33=\"),\n    @(\"78-17=\", \"39-39=\"),\n    @(\"70-52=\", \"79+17=\"),\n    @(\"81-49=\", \"12+5=\"),\n    @(\"3+73=\", \"92-71=\"),\n    @(\"96-85=\", \"13+80=\"),\n    @(\"54-48=\", \"1+59=\"),\n    @(\"31+56=\", \"87-55=\"),\n    @(\"23+28=\", \"99-20=\"),\n    @(\"27+33=\", \"44+43=\"),\n    @(\"61+31=\", \"66-9=\"),\n    @(\"18+18=\", \"67-64=\"),\n    @(\"0+66=\", \"91-59=\"),\n    @(\"90-71=\", \"32+3=\")\n)\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n$count = $paras.Count\n\n$idx = 0\n$applied = 0\n$mismatches = @()\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $paras.Item($i)\n    $r = $p.Range\n    $raw = $r.Text\n    $trimmed = $raw -replace \"[\\r\\a]\", \"\"\n    if ($trimmed.Length -eq 0) {\n        continue\n    }\n\n    if ($idx -ge $Replacements.Count) {\n        break\n    }\n\n    $pair = $Replacements[$idx]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $idx++\n\n    if ($trimmed -eq $oldText) {\n        $editRange = $p.Range\n        [void]$editRange.MoveEnd(1, -1)\n        $editRange.Text = $newText\n        $applied++\n    } else {\n        $mismatches += \"$i`: expected [$oldText] got [$trimmed]\"\n    }\n}\n\n\"applied=$applied total=$($Replacements.Count) mismatches=$($mismatches.Count)\"\nif ($mismatches.Count -gt 0) {\n    $mismatches -join \"`n\"\n}\n"}
